# "Coloquei os CSV dos testes" -- fill in the timing-test results (H:K and
# N columns) on the "Sheet2" tab (the sheet whose used range is A1:N7),
# matching the CSV the author pasted in.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Accuracy columns (H,I,J,K) for rows 4-7 ---------------------------
$ws2.Range("H4").Value = 0.96199999999999997
$ws2.Range("I4").Value = 0.94699999999999995
$ws2.Range("J4").Value = 0.97699999999999998
$ws2.Range("K4").Value = 0.96199999999999997

$ws2.Range("H5").Value = 0.97099999999999997
$ws2.Range("I5").Value = 0.96099999999999997
$ws2.Range("J5").Value = 0.98
$ws2.Range("K5").Value = 0.97099999999999997

$ws2.Range("H6").Value = 0.98599999999999999
$ws2.Range("I6").Value = 0.99
$ws2.Range("J6").Value = 0.98099999999999998
$ws2.Range("K6").Value = 0.98599999999999999

$ws2.Range("H7").Value = 0.98899999999999999
$ws2.Range("I7").Value = 0.99299999999999999
$ws2.Range("J7").Value = 0.98499999999999999
$ws2.Range("K7").Value = 0.98899999999999999

# --- Total-time-to-train column (N) for rows 2-7 -----------------------
$ws2.Range("N2").Value = "6s"
$ws2.Range("N3").Value = "6s"
$ws2.Range("N4").Value = "6s"
$ws2.Range("N5").Value = "1m50"
$ws2.Range("N6").Value = "1m43"
$ws2.Range("N7").Value = "1m42"

# Give the new N-column data cells the same look as the rest of the
# table (copy the format used by the other data cells in the row, which
# carries the fill/alignment), then drop the inner separators so only the
# left/right edges of the column are bordered -- same treatment as the
# header cell above them.
$fmtSource = $ws2.Range("H2")
$fmtSource.Copy()
$ws2.Range("N2:N7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($r = 2; $r -le 7; $r++) {
    $cell = $ws2.Cells.Item($r, 14)
    $cell.Borders.Item(8).LineStyle = -4142
    $cell.Borders.Item(9).LineStyle = -4142
}

# --- View bookkeeping ----------------------------------------------------
# Scroll "Sheet3" back to the top-left corner (it had drifted to A4).
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# Re-select the cell the author left active on "Sheet2" and make that the
# active tab again.
$ws2.Activate()
$ws2.Range("N8").Select()
